$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells (I1, J1) - new values with same header style as existing headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy formatting from an existing header cell (H1) onto the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data cells (plain numbers, no special style)
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 5

$ws.Range("I3").Value = 3
$ws.Range("J3").Value = 8

$ws.Range("I4").Value = 8
$ws.Range("J4").Value = 8
